$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# post/newsfeed (row 21) response description: add sender/recipient fields
$ws.Range("D21").Value = '{ "data" : [postId : {"postId":id, "text" : text, "url":"url", "timestamp" : timestamp,"senderId":id, "senderName":name, "senderLastname":lastname, "senderPicture":url, "senderUsername":username, "senderEmail":email,  "recipientId":id, "recipientName":name, "recipientLastname":lastname, "recipientPicture":url, "recipientUsername":username, "recipientEmail":email, "comments" : commentsNmbr, "likes": likesNmbr}, secondPost : {secondPost}, .. nthPost : {nthPost}] , "error" : [] }'

# post/getComments (row 22) response description: text->message, profilePicture->picture
$ws.Range("D22").Value = '{ "data" : [{"postId":id, "message" : text, "url":"url", "timestamp" : timestamp,"userId":id, "name":name, "lastName":lastname, "picture":url, "username":username}, {second comment}, …{nth comment}] , "error" : [] }'

# post/wall (row 23) response description: now matches the newsfeed (sender/recipient) format
$ws.Range("D23").Value = '{ "data" : [postId : {"postId":id, "text" : text, "url":"url", "timestamp" : timestamp,"senderId":id, "senderName":name, "senderLastname":lastname, "senderPicture":url, "senderUsername":username, "senderEmail":email,  "recipientId":id, "recipientName":name, "recipientLastname":lastname, "recipientPicture":url, "recipientUsername":username, "recipientEmail":email, "comments" : commentsNmbr, "likes": likesNmbr}, secondPost : {secondPost}, .. nthPost : {nthPost}] , "error" : [] }'

# post/wall and post/getComments become active
$ws.Range("F22").Value = "Active"
$ws.Range("F23").Value = "Active"

# Row heights adjust to fit the new wrapped text
$ws.Rows(21).RowHeight = 180
$ws.Rows(23).RowHeight = 180

# Update the view: selected cell in the frozen (bottom) pane moves to D22
$ws.Activate()
$ws.Range("D22").Select()
